$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(20, 0, 1, 0, 0, 0, 0, 0, 0, 0.05, 0, 0, 0, 0, 0, 1, 3, 0, 128, 0, 6, 0, 0, 0.7784000039100647, 0.7749999761581421),
    @(20, 0, 1, 0, 0, 0, 0, 0, 0, 0.05, 0, 0, 0, 0, 0, 1, 3, 0, 128, 0, 6, 0, 0, 0.7633000016212463, 0.7513999938964844),
    @(20, 1, 2.163987255, 0, 0, 0.102922781, 0.133724907, 0, 0, 0.010547912, 0, 0, 0, 0, 0, 0.087016579, 2.96260662, -0.111799513, 128, -0.001166751, 6, 0, 0.102167731, 0.8555999994277954, 0.8532000184059143),
    @(20, 0, 1, 0, 0, 0, 0, 0, 0, 0.05, 0, 0, 0, 0, 0, 1, 3, 0, 128, 0, 6, 0, 0, 0.7573999762535095, 0.7531999945640564),
    @(20, 1, 2.163987255, 0, 0, 0.102922781, 0.133724907, 0, 0, 0.010547912, 0, 0, 0, 0, 0, 0.087016579, 2.96260662, -0.111799513, 128, -0.001166751, 6, 0, 0.102167731, 0.8578000068664551, 0.843999981880188)
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
